$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A - shifts all existing columns right by one
$ws.Range("A1").EntireColumn.Insert()

# New column A header (row 3, the visible sub-header row) gets "Match ID"
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Data rows 4-14 get the Match ID value (15), bold styled like the header
for ($r = 4; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = 15
    $ws.Cells.Item($r, 1).Font.Bold = $true
}

# Row 15 is the hidden totals row - gets the value but keeps default styling
$ws.Range("A15").Value = 15
# Re-autofit the hidden row so writing to it doesn't leave a stray custom height
$ws.Rows.Item(15).AutoFit()

# Update the selection to match the new active range
$ws.Range("A3:A14").Select() | Out-Null
